# Apply the "cryptos list" update as described in the commit diff.
# Updates Price (D) and Volume(1h) (E) figures for most rows, and swaps
# the Cosmos/Toncoin (rows 28-29) and Stellar/EnergySwap (rows 40-41)
# entries (Coin name, Link, Price, Volume) to reflect their new ranking order.
# Price cells are plain text in the workbook (e.g. "47.994.44", "109.90");
# force text format before assignment so Excel does not coerce numeric-looking
# strings (losing formatting like trailing zeros, e.g. "109.90" -> 109.9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.994.44"
$ws.Range("E2").Value = "  +1.15%  "

$ws.Range("D3").Value = "2.511.84"
$ws.Range("E3").Value = "  +0.60%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("E5").Value = "  +0.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.90"
$ws.Range("E6").Value = "  +2.19%  "

$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.555"
$ws.Range("E9").Value = "  +1.96%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.91"
$ws.Range("E10").Value = "  +6.71%  "

$ws.Range("E11").Value = "  +0.37%  "

$ws.Range("E12").Value = "  +0.60%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.73"
$ws.Range("E13").Value = "  +1.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.26"
$ws.Range("E14").Value = "  +0.82%  "

$ws.Range("D15").Value = "2.905.58"
$ws.Range("E15").Value = "  +0.87%  "

$ws.Range("D16").Value = "2.510.54"
$ws.Range("E16").Value = "  +1.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.858"
$ws.Range("E17").Value = "  +1.21%  "

$ws.Range("D18").Value = "47.909.16"
$ws.Range("E18").Value = "  +1.20%  "

$ws.Range("E19").Value = "  +4.18%  "

$ws.Range("E20").Value = "  +1.01%  "

$ws.Range("E21").Value = "  +15.57%  "

$ws.Range("E22").Value = "  +0.68%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.96"
$ws.Range("E23").Value = "  +0.15%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "248.12"
$ws.Range("E24").Value = "  -1.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.57"
$ws.Range("E25").Value = "  -0.47%  "

$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("E27").Value = "  -1.00%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.31"
$ws.Range("E28").Value = "  +4.37%  "

$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.06"
$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("E30").Value = "  +2.96%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.30"
$ws.Range("E31").Value = "  -0.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.17"
$ws.Range("E33").Value = "  +1.69%  "

$ws.Range("E34").Value = "  -0.56%  "

$ws.Range("E35").Value = "  +0.21%  "

$ws.Range("E36").Value = "  +0.14%  "

$ws.Range("E37").Value = "  -0.20%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.68"
$ws.Range("E38").Value = "  +0.71%  "

$ws.Range("E39").Value = "  -0.12%  "

$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.112"
$ws.Range("E40").Value = "  +0.14%  "

$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.63"
$ws.Range("E41").Value = "  +6.08%  "

$ws.Range("E42").Value = "  -0.88%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "119.58"
$ws.Range("E43").Value = "  -2.08%  "

$ws.Range("E44").Value = "  +0.56%  "

$ws.Range("D45").Value = "2.004.82"
$ws.Range("E45").Value = "  +1.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.08"
$ws.Range("E46").Value = "  +2.53%  "

$ws.Range("E47").Value = "  -3.51%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.83"
$ws.Range("E48").Value = "  +1.59%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.06"
$ws.Range("E49").Value = "  -1.06%  "

$ws.Range("E50").Value = "  -0.99%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "57.19"
$ws.Range("E51").Value = "  +3.83%  "

